# "believe I've identified all bad zombies"
# Remove the bad zombie entry (sapling.id = 165) from the list entirely, and
# move the entry for sapling.id = 179 down to the bottom of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A before any edits.
$lastRow = $ws.Cells($ws.Rows.Count, 1).End(-4162).Row   # xlUp = -4162

# Locate the row holding the value 165 (bad zombie to remove) and the row
# holding the value 179 (needs to move to the end of the list).
$row165 = $null
$row179 = $null
for ($r = 1; $r -le $lastRow; $r++) {
    $val = $ws.Cells($r, 1).Value2
    if ($val -eq 165) { $row165 = $r }
    if ($val -eq 179) { $row179 = $r }
}

# Delete the entire row containing 165 - everything below shifts up one row.
$ws.Range("A" + $row165).EntireRow.Delete()

# Recompute where 179 now lives (it shifted up by one row since it was below 165).
if ($row179 -gt $row165) {
    $row179 = $row179 - 1
}
$lastRow = $lastRow - 1

# Move the 179 entry to the bottom of the list: write 179 into the row right
# after the current last row, then delete its old row (shifting the rest up).
$ws.Cells($lastRow + 1, 1).Value = 179
$ws.Range("A" + $row179).EntireRow.Delete()
